# ManageEmployees.xlsx update:
#   - Add two new data-provider sheets: FilterEmployees, AddEmployee
#   - Update the SearchForEmployees selection
#   - Make AddEmployee the active sheet (was SortEmployeeTableColumns)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. New sheet: FilterEmployees (inserted after ShowEmployeesPerPage)
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$filterSheet = $wb.Worksheets.Add($null, $lastSheet)
$filterSheet.Name = "FilterEmployees"

$filterSheet.Range("A1").Value = "TestScenario"
$filterSheet.Range("B1").Value = "Role"
$filterSheet.Range("A2").Value = "Filter Employees on Employee Page as Company Admin"
$filterSheet.Range("B2").Value = "COMPANY_ADMIN"
$filterSheet.Range("A1:B2").Select() | Out-Null

# ------------------------------------------------------------------
# 2. New sheet: AddEmployee (inserted after FilterEmployees)
# ------------------------------------------------------------------
$addSheet = $wb.Worksheets.Add($null, $filterSheet)
$addSheet.Name = "AddEmployee"

$addSheet.Range("A1").Value = "TestScenario"
$addSheet.Range("B1").Value = "Role"
$addSheet.Range("A2").Value = "Add New Employee as Company Admin"
$addSheet.Range("B2").Value = "COMPANY_ADMIN"

# ------------------------------------------------------------------
# 3. Tweak selection on SearchForEmployees
# ------------------------------------------------------------------
$searchSheet = $wb.Worksheets.Item("SearchForEmployees")
$searchSheet.Range("A1:B5").Select() | Out-Null

# ------------------------------------------------------------------
# 4. AddEmployee becomes the active tab (workbook activeTab -> 4)
# ------------------------------------------------------------------
$addSheet.Activate()
